$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.984.44'
$ws.Range('E2').Value = '  +1.03%  '

$ws.Range('D3').Value = '2.735.66'
$ws.Range('E3').Value = '  +3.42%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.57'
$ws.Range('E5').Value = '  +1.44%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.51'
$ws.Range('E6').Value = '  +6.34%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('E8').Value = '  +0.92%  '

$ws.Range('D9').Value = '2.736.15'
$ws.Range('E9').Value = '  +3.47%  '

$ws.Range('E10').Value = '  +3.57%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.371'
$ws.Range('E11').Value = '  +5.62%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.36'
$ws.Range('E12').Value = '  +1.84%  '

$ws.Range('E13').Value = '  -0.24%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.85'
$ws.Range('E14').Value = '  +3.19%  '

$ws.Range('D15').Value = '3.238.34'
$ws.Range('E15').Value = '  +3.51%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000192'
$ws.Range('E16').Value = '  +2.80%  '

$ws.Range('D17').Value = '68.939.48'
$ws.Range('E17').Value = '  +1.17%  '

$ws.Range('D18').Value = '2.717.24'
$ws.Range('E18').Value = '  +2.42%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.96'
$ws.Range('E19').Value = '  +5.26%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '375.29'
$ws.Range('E20').Value = '  +4.30%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.70'
$ws.Range('E21').Value = '  +4.92%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.56'
$ws.Range('E22').Value = '  +3.31%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.00'
$ws.Range('E23').Value = '  +5.69%  '

$ws.Range('E24').Value = '  +3.59%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.98'
$ws.Range('E25').Value = '  -1.28%  '

$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.13'
$ws.Range('E27').Value = '  +3.96%  '

$ws.Range('D28').Value = '2.869.78'
$ws.Range('E28').Value = '  +2.11%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000107'
$ws.Range('E29').Value = '  +3.60%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '587.82'
$ws.Range('E30').Value = '  +5.42%  '

$ws.Range('E31').Value = '  +0.44%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.38'
$ws.Range('E32').Value = '  +5.00%  '

$ws.Range('E33').Value = '  +5.56%  '

$ws.Range('E34').Value = '  +5.58%  '

$ws.Range('E35').Value = '  +4.58%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.64'
$ws.Range('E36').Value = '  +4.32%  '

$ws.Range('E37').Value = '  -0.01%  '

$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.51'
$ws.Range('E38').Value = '  +1.95%  '

$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '20.06'
$ws.Range('E39').Value = '  +1.54%  '

$ws.Range('E40').Value = '  +3.31%  '

$ws.Range('E41').Value = '  +4.09%  '

$ws.Range('E42').Value = '  +3.25%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.69'
$ws.Range('E43').Value = '  +3.27%  '

$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.00'
$ws.Range('E44').Value = '  +1.16%  '

$ws.Range('E45').Value = '  +0.07%  '

$ws.Range('D46').Value = '0.0₆0313'
$ws.Range('E46').Value = '  -1.98%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '41.07'
$ws.Range('E47').Value = '  +1.85%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '156.53'
$ws.Range('E48').Value = '  -0.10%  '

$ws.Range('E49').Value = '  +4.85%  '

$ws.Range('E50').Value = '  +7.30%  '

$ws.Range('E51').Value = '  +6.96%  '
